$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.117.46"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "3.749.29"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.80"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("D7").Value = "3.747.97"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.37"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.90"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("E14").Value = "  +1.41%  "
$ws.Range("D15").Value = "4.372.92"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").Value = "3.747.80"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").Value = "69.108.19"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("E18").Value = "  +1.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.43"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.114"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.04"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +9.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "492.43"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.728"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("E24").Value = "  +7.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.86"
$ws.Range("D25").ClearFormats()
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("E31").Value = "  +1.88%  "
$ws.Range("E32").Value = "  +1.87%  "
$ws.Range("D33").Value = "3.893.72"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.41"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.03%  "
$ws.Range("D35").Value = "3.680.58"
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.02"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.95"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.53%  "
$ws.Range("E40").Value = "  +3.22%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  +5.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.81"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "425.38"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.45"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.10"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.63"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("D50").Value = "2.787.77"
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("E51").Value = "  +0.14%  "
